# Upload new version with timestamp
# ---------------------------------------------------------------------------
# The sheet holds a "نواقص الأصناف" (missing/low-stock items) report built of
# a fixed-layout template: a header block (rows 1-6), one row per item
# (originally rows 7-14, 8 items), a totals row, and a footer row with the
# generation timestamp.
#
# The new version of the report has 4 more items (12 total), so 4 extra
# template rows need to be inserted right above the old totals row before
# the data / totals / footer are (re)written.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Make room: insert 4 new rows just above the current totals row (15),
#    cloning row 14's formatting/merges (Insert() alone only shifts cells,
#    it doesn't copy the per-row look, and a bare insert leaves blank
#    unmerged cells with default row height).
# ---------------------------------------------------------------------------
for ($i = 0; $i -lt 4; $i++) {
    $ws.Rows("15:15").Insert()
    $ws.Rows("14:14").Copy()
    $ws.Rows("15:15").PasteSpecial()
}

# Row heights exactly as they appear in the refreshed report.
$ws.Rows("7:7").RowHeight = 25.5
$ws.Rows("8:8").RowHeight = 24.75
$ws.Rows("9:9").RowHeight = 25.5
$ws.Rows("10:10").RowHeight = 24.75
$ws.Rows("11:11").RowHeight = 25.5
$ws.Rows("12:12").RowHeight = 25.5
$ws.Rows("13:13").RowHeight = 24.75
$ws.Rows("14:14").RowHeight = 25.5
$ws.Rows("15:15").RowHeight = 24.75
$ws.Rows("16:16").RowHeight = 25.5
$ws.Rows("17:17").RowHeight = 25.5
$ws.Rows("18:18").RowHeight = 24.75
$ws.Rows("19:19").RowHeight = 25.5
$ws.Rows("20:20").RowHeight = 16.5

# ---------------------------------------------------------------------------
# 2) Write the 12 data rows (7-18). Columns A (serial #) and the totals cell
#    are genuine numbers; every other data column (H/L/N/P/Q) is stored as
#    TEXT in the source report even when it looks numeric (e.g. "1", "46.00"),
#    so values are entered with a leading apostrophe to stop Excel from
#    re-typing them as numbers - this keeps each column's existing number
#    format (it was copied along with the row) instead of resetting it.
# ---------------------------------------------------------------------------
$rows = @(
    @{ Row = 7;  Serial = 1;  Name = "CORASORE 150MG 20 TAB";         Stock = "1:0";    Limit = "1"; Price = "46.00";  Sale = "46.0000"; Trans = "1:0" },
    @{ Row = 8;  Serial = 2;  Name = "DOLIPRANE 1 GM 15 TABS.";       Stock = "12:0";   Limit = "1"; Price = "48.00";  Sale = "48.0000"; Trans = "1:0" },
    @{ Row = 9;  Serial = 3;  Name = "EREC 100MG 12 F.C. TABLETS";    Stock = "1:10";   Limit = "1"; Price = "144.00"; Sale = "36.0000"; Trans = "0:3" },
    @{ Row = 10; Serial = 4;  Name = "FAWAR FRUIT 6 SACHETS";         Stock = "5:1";    Limit = "1"; Price = "24.00";  Sale = "7.9200";  Trans = "0:2" },
    @{ Row = 11; Serial = 5;  Name = "WATER FOR INJECTION AMP. 5 ML"; Stock = "8287:0"; Limit = "1"; Price = "2.00";   Sale = "2.0000";  Trans = "1:0" },
    @{ Row = 12; Serial = 6;  Name = "امواس لورد";                    Stock = "26:0";   Limit = "0"; Price = "15.00";  Sale = "15.0000"; Trans = "1:0" },
    @{ Row = 13; Serial = 7;  Name = "ببرونه صغير الجو";              Stock = "7:0";    Limit = "0"; Price = "20.00";  Sale = "20.0000"; Trans = "1:0" },
    @{ Row = 14; Serial = 8;  Name = "بلاستر مترسيلك 2 سم";           Stock = "32:0";   Limit = "0"; Price = "15.00";  Sale = "15.0000"; Trans = "1:0" },
    @{ Row = 15; Serial = 9;  Name = "ريكسونا حريمي بليه";            Stock = "6:0";    Limit = "0"; Price = "27.00";  Sale = "27.0000"; Trans = "1:0" },
    @{ Row = 16; Serial = 10; Name = "كالونا ";                       Stock = "0:0";    Limit = "0"; Price = "15.00";  Sale = "15.0000"; Trans = "1:0" },
    @{ Row = 17; Serial = 11; Name = "ماكينه حلاقه جليت فليكتور";     Stock = "14:0";   Limit = "0"; Price = "15.00";  Sale = "15.0000"; Trans = "1:0" },
    @{ Row = 18; Serial = 12; Name = "محلول جلوكوز 5%";               Stock = "20:0";   Limit = "0"; Price = "27.00";  Sale = "27.0000"; Trans = "1:0" }
)

foreach ($r in $rows) {
    $n = $r.Row
    $ws.Range("A$n").Value = $r.Serial
    $ws.Range("C$n").Value = "'" + $r.Name
    $ws.Range("H$n").Value = "'" + $r.Stock
    $ws.Range("L$n").Value = "'" + $r.Limit
    $ws.Range("N$n").Value = "'" + $r.Price
    $ws.Range("P$n").Value = "'" + $r.Sale
    $ws.Range("Q$n").Value = "'" + $r.Trans
}

# ---------------------------------------------------------------------------
# 3) Totals row (now row 19) and footer row (now row 20).
# ---------------------------------------------------------------------------
$ws.Range("P19").Value = 273.92

$ws.Range("A20").Value = "Friday, 12 September, 2025 12:40 PM"
